$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1998.3334
$ws.Range("I12").Value = 1998.3334
$ws.Range("K12").Value = 1998.3334
$ws.Range("M12").Value = -1828.3334
$ws.Range("H17").Value = 773756.9399999999
$ws.Range("J17").Value = 773756.9399999999
$ws.Range("L17").Value = 2321270.82
$ws.Range("N17").Value = -2321606.82
$ws.Range("H64").Value = 14428.571
$ws.Range("J64").Value = 16166.667
$ws.Range("L64").Value = 16166.667
$ws.Range("N64").Value = -16662.667
$ws.Range("H67").Value = 14428.571
$ws.Range("J67").Value = 16166.667
$ws.Range("L67").Value = 16166.667
$ws.Range("N67").Value = -17882.667
$ws.Range("H96").Value = 1983.8182
$ws.Range("I96").Value = 2007
$ws.Range("K96").Value = 6021
$ws.Range("M96").Value = -4648
$ws.Range("H137").Value = 1392.8889
$ws.Range("I137").Value = 1087.5
$ws.Range("J137").Value = 1637.2
$ws.Range("K137").Value = 3262.5
$ws.Range("L137").Value = 4911.6
$ws.Range("M137").Value = -712.5
$ws.Range("N137").Value = -10011.6
$ws.Range("H138").Value = 3632.1265
$ws.Range("J138").Value = 4109.3223
$ws.Range("L138").Value = 12327.9669
$ws.Range("N138").Value = -22607.9669
$ws.Range("H141").Value = 1666.6666
$ws.Range("I141").Value = 1666.6666
$ws.Range("K141").Value = 4999.9998
$ws.Range("M141").Value = 180.0002000000004

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4035.42
$ws.Range("I32").Value = 3527.761
$ws.Range("K32").Value = 3527.761
$ws.Range("M32").Value = -3240.761
$ws.Range("H63").Value = 2600
$ws.Range("I63").Value = 2700
$ws.Range("K63").Value = 2700
$ws.Range("M63").Value = -2014
$ws.Range("H66").Value = 2600
$ws.Range("I66").Value = 2700
$ws.Range("K66").Value = 13500
$ws.Range("M66").Value = -10068
$ws.Range("H122").Value = 4442.3335
$ws.Range("I122").Value = 4478.778
$ws.Range("J122").Value = 4333
$ws.Range("K122").Value = 13436.334
$ws.Range("L122").Value = 12999
$ws.Range("M122").Value = -10986.334
$ws.Range("N122").Value = -17899

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4097.4814
$ws.Range("I20").Value = 3578.75
$ws.Range("J20").Value = 4512.467
$ws.Range("K20").Value = 3578.75
$ws.Range("L20").Value = 4512.467
$ws.Range("M20").Value = -3331.75
$ws.Range("N20").Value = -5006.467
$ws.Range("H86").Value = 2000
$ws.Range("I86").Value = 2000
$ws.Range("J86").Value = 2000
$ws.Range("K86").Value = 2000
$ws.Range("L86").Value = 2000
$ws.Range("M86").Value = -877
$ws.Range("N86").Value = -4246
$ws.Range("H89").Value = 2000
$ws.Range("I89").Value = 2000
$ws.Range("J89").Value = 2000
$ws.Range("K89").Value = 10000
$ws.Range("L89").Value = 10000
$ws.Range("M89").Value = -4384
$ws.Range("N89").Value = -21232

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 504.42856
$ws.Range("I11").Value = 225
$ws.Range("J11").Value = 1203
$ws.Range("K11").Value = 225
$ws.Range("L11").Value = 1203
$ws.Range("M11").Value = -85
$ws.Range("N11").Value = -1483
$ws.Range("H12").Value = 1035.174
$ws.Range("I12").Value = 433.375
$ws.Range("J12").Value = 2410.7144
$ws.Range("K12").Value = 433.375
$ws.Range("L12").Value = 2410.7144
$ws.Range("M12").Value = -263.375
$ws.Range("N12").Value = -2750.7144
$ws.Range("H14").Value = 1217.25
$ws.Range("I14").Value = 648.1667
$ws.Range("J14").Value = 2924.5
$ws.Range("K14").Value = 648.1667
$ws.Range("L14").Value = 2924.5
$ws.Range("M14").Value = -478.1667
$ws.Range("N14").Value = -3264.5
$ws.Range("H31").Value = 26822.219
$ws.Range("I31").Value = 36334.242
$ws.Range("K31").Value = 36334.242
$ws.Range("M31").Value = -36039.242
$ws.Range("H34").Value = 26822.219
$ws.Range("I34").Value = 36334.242
$ws.Range("K34").Value = 36334.242
$ws.Range("M34").Value = -36132.242
$ws.Range("H56").Value = 21333.334
$ws.Range("I56").Value = 17500
$ws.Range("J56").Value = 29000
$ws.Range("K56").Value = 17500
$ws.Range("L56").Value = 29000
$ws.Range("M56").Value = -16655
$ws.Range("N56").Value = -30690
$ws.Range("H59").Value = 29000
$ws.Range("J59").Value = 29000
$ws.Range("L59").Value = 29000
$ws.Range("N59").Value = -31290
$ws.Range("H60").Value = 12984.714
$ws.Range("I60").Value = 7724
$ws.Range("K60").Value = 7724
$ws.Range("M60").Value = -7213

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 37787748
$ws.Range("I4").Value = 40433300
$ws.Range("K4").Value = 121299900
$ws.Range("M4").Value = -121299788
$ws.Range("H11").Value = 107.14286
$ws.Range("I11").Value = 93.75
$ws.Range("J11").Value = 150
$ws.Range("K11").Value = 281.25
$ws.Range("L11").Value = 450
$ws.Range("M11").Value = -141.25
$ws.Range("N11").Value = -730
$ws.Range("H122").Value = 994.4074000000001
$ws.Range("I122").Value = 1280.6154
$ws.Range("J122").Value = 728.6429000000001
$ws.Range("K122").Value = 11525.5386
$ws.Range("L122").Value = 6557.7861
$ws.Range("M122").Value = -9075.5386
$ws.Range("N122").Value = -11457.7861
$ws.Range("H131").Value = 22945.084
$ws.Range("I131").Value = 125658.25
$ws.Range("J131").Value = 2402.45
$ws.Range("K131").Value = 376974.75
$ws.Range("L131").Value = 7207.349999999999
$ws.Range("M131").Value = -371934.75
$ws.Range("N131").Value = -17287.35

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 15750
$ws.Range("J15").Value = 15750
$ws.Range("L15").Value = 15750
$ws.Range("N15").Value = -16326
$ws.Range("H57").Value = 8003.923
$ws.Range("J57").Value = 8499.5
$ws.Range("L57").Value = 8499.5
$ws.Range("N57").Value = -10139.5
$ws.Range("H70").Value = 8088.037
$ws.Range("I70").Value = 5423.091
$ws.Range("J70").Value = 19813.8
$ws.Range("K70").Value = 5423.091
$ws.Range("L70").Value = 19813.8
$ws.Range("M70").Value = -5153.091
$ws.Range("N70").Value = -20353.8
$ws.Range("H73").Value = 8088.037
$ws.Range("I73").Value = 5423.091
$ws.Range("J73").Value = 19813.8
$ws.Range("K73").Value = 5423.091
$ws.Range("L73").Value = 19813.8
$ws.Range("M73").Value = -4487.091
$ws.Range("N73").Value = -21685.8
$ws.Range("H80").Value = 2690.375
$ws.Range("J80").Value = 2398.5
$ws.Range("L80").Value = 2398.5
$ws.Range("N80").Value = -4394.5
$ws.Range("H81").Value = 15750
$ws.Range("J81").Value = 15750
$ws.Range("L81").Value = 15750
$ws.Range("N81").Value = -17746
$ws.Range("H83").Value = 2690.375
$ws.Range("J83").Value = 2398.5
$ws.Range("L83").Value = 11992.5
$ws.Range("N83").Value = -21976.5
$ws.Range("H84").Value = 15750
$ws.Range("J84").Value = 15750
$ws.Range("L84").Value = 47250
$ws.Range("N84").Value = -57234
$ws.Range("H102").Value = 37038296
$ws.Range("I102").Value = 963.1905
$ws.Range("K102").Value = 963.1905
$ws.Range("M102").Value = 658.8095
$ws.Range("H135").Value = 49666.223
$ws.Range("J135").Value = 49666.223
$ws.Range("L135").Value = 49666.223
$ws.Range("N135").Value = -59806.223

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4718.5
$ws.Range("I61").Value = 4672.5293
$ws.Range("J61").Value = 5500
$ws.Range("K61").Value = 4672.5293
$ws.Range("L61").Value = 5500
$ws.Range("M61").Value = -4470.5293
$ws.Range("N61").Value = -5904
$ws.Range("H113").Value = 4718.5
$ws.Range("I113").Value = 4672.5293
$ws.Range("J113").Value = 5500
$ws.Range("K113").Value = 4672.5293
$ws.Range("L113").Value = 5500
$ws.Range("M113").Value = -2502.5293
$ws.Range("N113").Value = -9840
$ws.Range("H122").Value = 6017.8486
$ws.Range("I122").Value = 5132.227
$ws.Range("J122").Value = 7789.091
$ws.Range("K122").Value = 15396.681
$ws.Range("L122").Value = 23367.273
$ws.Range("M122").Value = -12946.681
$ws.Range("N122").Value = -28267.273
$ws.Range("H132").Value = 4541.684
$ws.Range("I132").Value = 4517.9375
$ws.Range("J132").Value = 4668.3335
$ws.Range("K132").Value = 13553.8125
$ws.Range("L132").Value = 14005.0005
$ws.Range("M132").Value = -11023.8125
$ws.Range("N132").Value = -19065.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 4568
$ws.Range("I7").Value = 4568
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 4568
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -4455
$ws.Range("N7").ClearContents()
$ws.Range("H54").Value = 17450
$ws.Range("J54").Value = 17450
$ws.Range("L54").Value = 17450
$ws.Range("N54").Value = -18490
$ws.Range("H95").Value = 58516.082
$ws.Range("J95").Value = 58516.082
$ws.Range("L95").Value = 58516.082
$ws.Range("N95").Value = -64008.082
$ws.Range("H122").Value = 2958.2058
$ws.Range("I122").Value = 2867.1785
$ws.Range("K122").Value = 8601.5355
$ws.Range("M122").Value = -6151.5355
